# Generate Report for Handoff
# Adds a new handoff-status row (for file "c738f36e-d1e0-4800-adfe-ef5cbad142d7")
# to the Overview, zh-cn and de-de worksheets of the localization-status report.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # packed R+G*256+B*65536 for RGB(0x64,0x95,0xED) -> renders as FF6495ED

$newFileId = "c738f36e-d1e0-4800-adfe-ef5cbad142d7"
$newMdName = "$newFileId.md"
$newHash   = "7bb04417058798f7cabaf1f5017023dbdbf327cc"
$zhXlfName = "$newFileId.$newHash.zh-cn.xlf"
$deXlfName = "$newFileId.$newHash.de-de.xlf"

$mdCommit = "8186d62ccf7150a9fe64afdde09ed3cc954bd983"
$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$newMdName"

$zhCommit = "795c348b1d284dc2af595b278676a86b8fa9ea01"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"

$deCommit = "e2163d9acd0f504aa293c244c9d22458b0992282"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-28-20 00:28:55"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $newMdName)
$wsOverview.Range("A3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newMdName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = $zhXlfName
$wsZh.Range("E3").Value = "2016-03-20 00:28:52"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $newMdName)
$wsZh.Range("A3").Font.Color = $hyperlinkColor

$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl, [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Range("B3").Font.Color = $hyperlinkColor

$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$wsZh.Range("D3").Font.Color = $hyperlinkColor

$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newMdName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = $deXlfName
$wsDe.Range("E3").Value = "2016-03-20 00:28:55"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $newMdName)
$wsDe.Range("A3").Font.Color = $hyperlinkColor

$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl, [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Range("B3").Font.Color = $hyperlinkColor

$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$wsDe.Range("D3").Font.Color = $hyperlinkColor

$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Report rows appended for $newFileId"
